$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F21").Value = '93_referral_statement'
$ws.Range("F23").Value = '18_hazards_to_humans_and_domestic_animals'
$ws.Range("F25").Value = 'ppe'
$ws.Range("F29").Value = 'application instructions || env warning - species'
$ws.Range("F30").Value = 'env warning - water'
$ws.Range("F31").Value = 'env warning - water || off target movement'
$ws.Range("F32").Value = 'safety procedures'
$ws.Range("F33").Value = 'use restrictions'
$ws.Range("F34").Value = '32_physical_and_chemical_hazards'
$ws.Range("F38").Value = 'application instructions'
$ws.Range("F39").Value = 'application instructions'
$ws.Range("F40").Value = '135_product_information'
$ws.Range("F45").Value = 'use restrictions'
$ws.Range("F49").Value = 'application instructions'
$ws.Range("F50").Value = 'application instructions'
$ws.Range("F52").Value = 'application instructions'
$ws.Range("F53").Value = 'application instructions'
$ws.Range("F54").Value = 'application instructions'
$ws.Range("F56").Value = 'use restrictions || off target movement'
$ws.Range("F57").Value = 'off target movement'
$ws.Range("F58").Value = 'off target movement'
$ws.Range("F59").Value = 'off target movement'
$ws.Range("F60").Value = 'off target movement'
$ws.Range("F61").Value = 'off target movement'
$ws.Range("F62").Value = '172_sensitive_areas'
$ws.Range("F63").Value = 'mixing'
$ws.Range("F65").Value = 'irrigation'
$ws.Range("F66").Value = 'application instructions'
$ws.Range("F67").Value = 'application instructions'
$ws.Range("F70").Value = 'mixing'
$ws.Range("F72").Value = 'mixing'
$ws.Range("F73").Value = 'safety procedures'
$ws.Range("F104").Value = 'use restrictions'
$ws.Range("F107").Value = 'application instructions'
$ws.Range("F108").Value = 'application instructions'
$ws.Range("F111").Value = 'application instructions'
$ws.Range("F113").Value = 'mixing'
$ws.Range("F117").Value = 'mixing'
$ws.Range("F123").Value = 'application instructions'
$ws.Range("F124").Value = 'application instructions'
$ws.Range("F126").Value = 'application instructions'
$ws.Range("F129").Value = '154_pesticide_storage'
